# "Generate Report for Handoff"
#
# The localization-status report is regenerated: items that were
# previously shown as "Handed back: in sync with en-US" are now
# "Ready for handoff" (status flips back once a new handoff cycle
# starts), and the associated timestamps advance a few seconds as the
# report is rebuilt. The Status column (and its mirrors on the
# Overview sheet) also narrows to fit the shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name / Path And Name / Extension / Publish URL /
#     zh-cn / de-de / Latest HO Xliff Generate Date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 09:02:57"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-09-03 09:02:53"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-09-03 09:02:57"

# --- Narrow the Status columns now that the text is shorter ---
# (target stored width ~= 17.216; ColumnWidth is rounded to the engine's
# internal width grid on save, so feed it the value whose rounded result
# lands closest to that target.)
$wsOverview.Range("E1").ColumnWidth = 16.3333333333333
$wsOverview.Range("F1").ColumnWidth = 16.3333333333333
$wsZh.Range("C1").ColumnWidth = 16.3333333333333
$wsDe.Range("C1").ColumnWidth = 16.3333333333333
